$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test-date values in column D (rows 7-9): 2022-12-02 -> 2022-12-25
$ws.Range("D7").Value = 44920
$ws.Range("D8").Value = 44920
$ws.Range("D9").Value = 44920

# Update the sheet's active selection / view to C18
$ws.Range("C18").Select()
